$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.01019667505792274
$ws.Range("C2").Value = 0.2651573259267139
$ws.Range("D2").Value = 0.09370580426362109
$ws.Range("E2").Value = 0.3061140380048277
$ws.Range("F2").Value = 0.3174932645123695

$ws.Range("B3").Value = 0.01598200330380756
$ws.Range("C3").Value = 0.2434613949423638
$ws.Range("D3").Value = 0.07230275474488276
$ws.Range("E3").Value = 0.2688917156494093
$ws.Range("F3").Value = 0.2829356627206967

$ws.Range("B4").Value = -0.02630075056802729
$ws.Range("C4").Value = 0.252529210783159
$ws.Range("D4").Value = 0.06708865167695137
$ws.Range("E4").Value = 0.2590147711559157
$ws.Range("F4").Value = 0.2822699180497485

$ws.Range("B5").Value = 0.009192281051027068
$ws.Range("C5").Value = 0.267502854063783
$ws.Range("D5").Value = 0.07164227496319067
$ws.Range("E5").Value = 0.2676607460259922
$ws.Range("F5").Value = 0.3783061641905128
